# Applies the cryptos-list price/volume refresh described by the commit diff.
# Rows 28/29 additionally swap Coin name + Link (InjectiveProtocol <-> WEMIXToken).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.797.10"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").Value = "2.239.38"
$ws.Range("E3").Value = "  +2.24%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'270.29"
$ws.Range("E5").Value = "  +4.46%  "

$ws.Range("D6").Value = "'92.36"
$ws.Range("E6").Value = "  +13.27%  "

$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  +1.12%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").Value = "'0.623"
$ws.Range("E9").Value = "  +4.98%  "

$ws.Range("D10").Value = "'46.21"
$ws.Range("E10").Value = "  +7.55%  "

$ws.Range("D11").Value = "'0.0956"
$ws.Range("E11").Value = "  +3.99%  "

$ws.Range("E12").Value = "  +18.68%  "

$ws.Range("E13").Value = "  +2.19%  "

$ws.Range("D14").Value = "2.573.47"
$ws.Range("E14").Value = "  +2.19%  "

$ws.Range("D15").Value = "'15.05"
$ws.Range("E15").Value = "  +5.19%  "

$ws.Range("D16").Value = "2.271.74"
$ws.Range("E16").Value = "  +3.73%  "

$ws.Range("D17").Value = "'0.806"
$ws.Range("E17").Value = "  +3.67%  "

$ws.Range("D18").Value = "43.801.29"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("E19").Value = "  +2.42%  "

$ws.Range("E20").Value = "  +2.96%  "

$ws.Range("D21").Value = "'70.79"
$ws.Range("E21").Value = "  +0.89%  "

$ws.Range("E22").Value = "  -2.64%  "

$ws.Range("D23").Value = "'233.94"
$ws.Range("E23").Value = "  +1.50%  "

$ws.Range("D24").Value = "'9.03"
$ws.Range("E24").Value = "  +0.66%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").Value = "'11.45"
$ws.Range("E26").Value = "  +7.53%  "

$ws.Range("D27").Value = "'2.50"
$ws.Range("E27").Value = "  +12.38%  "

$ws.Range("B28").Value = "WEMIXToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D28").Value = "'3.54"
$ws.Range("E28").Value = "  +5.43%  "

$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'41.51"
$ws.Range("E29").Value = "  -0.77%  "

$ws.Range("E30").Value = "  +1.11%  "

$ws.Range("D31").Value = "'172.44"
$ws.Range("E31").Value = "  -0.37%  "

$ws.Range("D32").Value = "'0.0920"
$ws.Range("E32").Value = "  +5.19%  "

$ws.Range("D33").Value = "'20.96"
$ws.Range("E33").Value = "  +2.86%  "

$ws.Range("D34").Value = "'5.51"
$ws.Range("E34").Value = "  +4.37%  "

$ws.Range("E35").Value = "  +1.51%  "

$ws.Range("E36").Value = "  +1.20%  "

$ws.Range("E37").Value = "  -0.33%  "

$ws.Range("D38").Value = "'4.31"
$ws.Range("E38").Value = "  -3.96%  "

$ws.Range("D39").Value = "'3.52"
$ws.Range("E39").Value = "  +25.09%  "

$ws.Range("D40").Value = "'0.231"
$ws.Range("E40").Value = "  +16.74%  "

$ws.Range("D41").Value = "'12.76"
$ws.Range("E41").Value = "  -3.11%  "

$ws.Range("E42").Value = "  +4.69%  "

$ws.Range("D43").Value = "'63.61"
$ws.Range("E43").Value = "  +1.39%  "

$ws.Range("D44").Value = "'5.37"
$ws.Range("E44").Value = "  -1.22%  "

$ws.Range("D45").Value = "'0.0996"
$ws.Range("E45").Value = "  +1.09%  "

$ws.Range("D46").Value = "'8.37"
$ws.Range("E46").Value = "  +1.90%  "

$ws.Range("D47").Value = "'100.34"
$ws.Range("E47").Value = "  -0.74%  "

$ws.Range("E48").Value = "  +4.42%  "

$ws.Range("E49").Value = "  +1.68%  "

$ws.Range("D50").Value = "'0.443"
$ws.Range("E50").Value = "  +1.45%  "

$ws.Range("D51").Value = "2.458.57"
$ws.Range("E51").Value = "  +2.12%  "
